$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Structural shape: insert the new "Time for 1000ul" column, drop the old
#    "no-blowouts"-duplicate + "difference" columns, and insert two new data
#    rows (viscosity standards 204 and 398) ahead of the existing 505 row.
# ---------------------------------------------------------------------------
$ws.Range("F1").EntireColumn.Insert()
$ws.Range("J1:O9").EntireColumn.Delete()
$ws.Rows.Item(4).Insert()
$ws.Rows.Item(4).Insert()

# Fix the merged header band: B2:E2 grows to B2:F2
$ws.Range("B2:E2").UnMerge()
$ws.Range("B2:F2").Merge()

# Restore the (now-empty) trailing merge placeholders that used to hold the
# "blowouts" / "difference" header bands.
$ws.Range("J2:L2").Merge()
$ws.Range("M2:O2").Merge()

Write-Host "structure done"

# ---------------------------------------------------------------------------
# 2. Header strings
# ---------------------------------------------------------------------------
$ws.Range("G2").Value = "% error (blowouts before)"

$ws.Range("A3").Value = "Viscosity standard:"
$ws.Range("B3").Value = "Aspiration rate:"
$ws.Range("C3").Value = "Dispense rate:"
$ws.Range("D3").Value = "Delay aspirate:"
$ws.Range("E3").Value = "Delay dispense:"
$ws.Range("F3").Value = "Time for 1000ul"
$ws.Range("G3").Value = "1000(ul)"
$ws.Range("H3").Value = "500(ul)"
$ws.Range("I3").Value = "300(ul)"

Write-Host "headers done"

# ---------------------------------------------------------------------------
# 3. Data rows
# ---------------------------------------------------------------------------
# Row 4 - viscosity standard 204
$ws.Range("A4").Value = 204.0
$ws.Range("B4").Value = 185.0
$ws.Range("C4").Value = 18.0
$ws.Range("D4").Value = 10.0
$ws.Range("E4").Value = 10.0
$ws.Range("F4").Value = 80.960961
$ws.Range("G4").Value = -0.992463
$ws.Range("H4").Value = -0.665249
$ws.Range("I4").Value = -3.022089

# Row 5 - viscosity standard 398
$ws.Range("A5").Value = 398.0
$ws.Range("B5").Value = 80.0
$ws.Range("C5").Value = 15.0
$ws.Range("D5").Value = 10.0
$ws.Range("E5").Value = 10.0
$ws.Range("F5").Value = 99.166667
$ws.Range("G5").Value = -1.767441
$ws.Range("H5").Value = -1.275127
$ws.Range("I5").Value = -3.49754

# Row 6 - viscosity standard 505 (was row 4)
$ws.Range("A6").Value = 505.0
$ws.Range("B6").Value = 20.0
$ws.Range("C6").Value = 5.0
$ws.Range("D6").Value = 10.0
$ws.Range("E6").Value = 10.0
$ws.Range("F6").Value = ""
$ws.Range("G6").Value = -2.215764
$ws.Range("H6").Value = -1.799589
$ws.Range("I6").Value = -4.330199

# Row 7 - viscosity standard 817 (was row 5)
$ws.Range("A7").Value = 817.0
$ws.Range("B7").Value = 10.0
$ws.Range("C7").Value = 5.0
$ws.Range("D7").Value = 10.0
$ws.Range("E7").Value = 10.0
$ws.Range("F7").Value = 320.0
$ws.Range("G7").Value = -2.759218
$ws.Range("H7").Value = -3.829765
$ws.Range("I7").Value = -4.086789

# Row 8 - viscosity standard 1275 (new)
$ws.Range("A8").Value = 1275.0
$ws.Range("B8").Value = 4.0
$ws.Range("C8").Value = 2.5
$ws.Range("D8").Value = 10.0
$ws.Range("E8").Value = 10.0
$ws.Range("F8").Value = 670.0
$ws.Range("G8").Value = -2.46707797254956
$ws.Range("H8").Value = -3.7964279316069
$ws.Range("I8").Value = -2.44205883906668

# Row 9 - trailing blank-ish note
$ws.Range("D9").Value = "`n"

# Match the "0.0" look for the aspiration/dispense-rate columns (B:C) on
# every data row.
$ws.Range("B4:C8").NumberFormat = "0.0"

Write-Host "data done"
